$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.383.98'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.335.69'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.04'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.19'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.581'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.334.50'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.64'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.751.26'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.313.23'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.329.80'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.68'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.11'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '315.21'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.22'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.171'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.94'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.28'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.33'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.08'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.13'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.10'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '324.48'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.15'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.76'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0950'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.38'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0500'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0220'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.93'

$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -1.32%  '
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("E11").Value = '  +1.80%  '
$ws.Range("E13").Value = '  +0.62%  '
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("E22").Value = '  -3.44%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("E25").Value = '  -1.19%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  +0.40%  '
$ws.Range("E28").Value = '  +1.14%  '
$ws.Range("E29").Value = '  +9.67%  '
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("E32").Value = '  +0.42%  '
$ws.Range("E33").Value = '  +2.14%  '
$ws.Range("E34").Value = '  -3.31%  '
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("E39").Value = '  -1.76%  '
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E44").Value = '  +1.49%  '
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("E48").Value = '  +0.24%  '
$ws.Range("E49").Value = '  +0.80%  '
$ws.Range("E50").Value = '  +4.55%  '
$ws.Range("E51").Value = '  -0.81%  '
